{"js": "// Update the date heading and the 25 \"three-digit \u00f7 one-digit\" problems\n// in the table, cell-by-cell (positional), so no text collisions occur\n// between old/new values.\n\nconst body = context.document.body;\n\n// --- 1) Date heading (first paragraph) ---------------------------------\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst titlePara = paras.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\n\nif (titlePara.text.trim() === \"2024-05-01 Wednesday\") {\n  titlePara.getRange().insertText(\"2024-05-02 Thursday\", \"Replace\");\n  await context.sync();\n}\n\n// --- 2) Table cells ------------------------------------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New values, row-major, only the 5 populated rows (others are blank\n// spacer rows in the grid) -- indices 0, 4, 8, 12, 16 of table.rows.\nconst newValuesByRow = {\n  0: [\"959\u00f79=\", \"650\u00f78=\", \"748\u00f74=\", \"939\u00f76=\", \"752\u00f75=\"],\n  4: [\"461\u00f76=\", \"180\u00f73=\", \"155\u00f77=\", \"981\u00f76=\", \"914\u00f72=\"],\n  8: [\"965\u00f78=\", \"346\u00f74=\", \"580\u00f79=\", \"609\u00f76=\", \"148\u00f74=\"],\n  12: [\"514\u00f72=\", \"627\u00f78=\", \"686\u00f79=\", \"338\u00f78=\", \"690\u00f74=\"],\n  16: [\"761\u00f73=\", \"718\u00f78=\", \"625\u00f75=\", \"116\u00f73=\", \"119\u00f72=\"],\n};\n\nfor (const rowIndexStr of Object.keys(newValuesByRow)) {\n  const rowIndex = Number(rowIndexStr);\n  const rowValues = newValuesByRow[rowIndex];\n  for (let col = 0; col < rowValues.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    cell.value = rowValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 \"three-digit / one-digit\" problems\n# in the table, cell-by-cell (positional), so no text collisions occur\n# between old/new values (e.g. one cell's new value equals another\n# cell's old value).\n\n$d = $word.ActiveDocument\n\n# --- 1) Date heading (first paragraph) -----------------------------------\n$p1 = $d.Paragraphs(1)\nif ($p1.Range.Text.Trim() -eq \"2024-05-01 Wednesday\") {\n    $p1.Range.Text = \"2024-05-02 Thursday\"\n}\n\n# --- 2) Table cells --------------------------------------------------------\n$t = $d.Tables.Item(1)\n\n$t.Cell(1,1).Range.Text = \"959\u00f79=\"\n$t.Cell(1,2).Range.Text = \"650\u00f78=\"\n$t.Cell(1,3).Range.Text = \"748\u00f74=\"\n$t.Cell(1,4).Range.Text = \"939\u00f76=\"\n$t.Cell(1,5).Range.Text = \"752\u00f75=\"\n\n$t.Cell(5,1).Range.Text = \"461\u00f76=\"\n$t.Cell(5,2).Range.Text = \"180\u00f73=\"\n$t.Cell(5,3).Range.Text = \"155\u00f77=\"\n$t.Cell(5,4).Range.Text = \"981\u00f76=\"\n$t.Cell(5,5).Range.Text = \"914\u00f72=\"\n\n$t.Cell(9,1).Range.Text = \"965\u00f78=\"\n$t.Cell(9,2).Range.Text = \"346\u00f74=\"\n$t.Cell(9,3).Range.Text = \"580\u00f79=\"\n$t.Cell(9,4).Range.Text = \"609\u00f76=\"\n$t.Cell(9,5).Range.Text = \"148\u00f74=\"\n\n$t.Cell(13,1).Range.Text = \"514\u00f72=\"\n$t.Cell(13,2).Range.Text = \"627\u00f78=\"\n$t.Cell(13,3).Range.Text = \"686\u00f79=\"\n$t.Cell(13,4).Range.Text = \"338\u00f78=\"\n$t.Cell(13,5).Range.Text = \"690\u00f74=\"\n\n$t.Cell(17,1).Range.Text = \"761\u00f73=\"\n$t.Cell(17,2).Range.Text = \"718\u00f78=\"\n$t.Cell(17,3).Range.Text = \"625\u00f75=\"\n$t.Cell(17,4).Range.Text = \"116\u00f73=\"\n$t.Cell(17,5).Range.Text = \"119\u00f72=\"\n"}
